$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 4125
$ws.Range("I18").Value = 5750
$ws.Range("J18").Value = 2500
$ws.Range("K18").Value = 5750
$ws.Range("L18").Value = 2500
$ws.Range("M18").Value = -5466
$ws.Range("N18").Value = -3068
$ws.Range("H33").Value = 422.57895
$ws.Range("I33").Value = 426.5
$ws.Range("J33").Value = 411.6
$ws.Range("K33").Value = 426.5
$ws.Range("L33").Value = 411.6
$ws.Range("M33").Value = -197.5
$ws.Range("N33").Value = -869.6
$ws.Range("H38").Value = 672.0909
$ws.Range("I38").Value = 154.77777
$ws.Range("J38").Value = 3000
$ws.Range("K38").Value = 464.33331
$ws.Range("L38").Value = 9000
$ws.Range("M38").Value = -92.33330999999998
$ws.Range("N38").Value = -9744
$ws.Range("H39").Value = 453.16666
$ws.Range("I39").Value = 243.8
$ws.Range("J39").Value = 1500
$ws.Range("K39").Value = 731.4000000000001
$ws.Range("L39").Value = 4500
$ws.Range("M39").Value = -435.4000000000001
$ws.Range("N39").Value = -5092
$ws.Range("H40").Value = 2305.2
$ws.Range("J40").Value = 2305.2
$ws.Range("L40").Value = 2305.2
$ws.Range("N40").Value = -2655.2
$ws.Range("H43").Value = 7943771
$ws.Range("I43").Value = 40001
$ws.Range("J43").Value = 9261066
$ws.Range("K43").Value = 40001
$ws.Range("L43").Value = 9261066
$ws.Range("M43").Value = -39932
$ws.Range("N43").Value = -9261204
$ws.Range("H106").Value = 3001.4614
$ws.Range("I106").Value = 3134.7778
$ws.Range("K106").Value = 3134.7778
$ws.Range("M106").Value = -2503.7778
$ws.Range("H137").Value = 1087.2222
$ws.Range("I137").Value = 1085.8572
$ws.Range("J137").Value = 1092
$ws.Range("K137").Value = 3257.5716
$ws.Range("L137").Value = 3276
$ws.Range("M137").Value = -707.5715999999998
$ws.Range("N137").Value = -8376
$ws.Range("H138").Value = 588809.4
$ws.Range("I138").Value = 715.5
$ws.Range("J138").Value = 1088689.1
$ws.Range("K138").Value = 2146.5
$ws.Range("L138").Value = 3266067.3
$ws.Range("M138").Value = 2993.5
$ws.Range("N138").Value = -3276347.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9256.167
$ws.Range("I2").Value = 894.2222
$ws.Range("K2").Value = 894.2222
$ws.Range("M2").Value = -781.2222
$ws.Range("H45").Value = 2021.1428
$ws.Range("I45").Value = 2021.1428
$ws.Range("K45").Value = 2021.1428
$ws.Range("M45").Value = -1644.1428
$ws.Range("H61").Value = 1685.4615
$ws.Range("I61").Value = 1339.7
$ws.Range("J61").Value = 2838
$ws.Range("K61").Value = 1339.7
$ws.Range("L61").Value = 2838
$ws.Range("M61").Value = -1127.7
$ws.Range("N61").Value = -3262
$ws.Range("H116").Value = 9256.167
$ws.Range("I116").Value = 894.2222
$ws.Range("K116").Value = 894.2222
$ws.Range("M116").Value = 1399.7778
$ws.Range("H122").Value = 1580
$ws.Range("I122").Value = 1423.3572
$ws.Range("K122").Value = 4270.071599999999
$ws.Range("M122").Value = -1820.071599999999
$ws.Range("H136").Value = 1685.4615
$ws.Range("I136").Value = 1339.7
$ws.Range("J136").Value = 2838
$ws.Range("K136").Value = 4019.1
$ws.Range("L136").Value = 8514
$ws.Range("M136").Value = -1469.1
$ws.Range("N136").Value = -13614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9256.167
$ws.Range("I3").Value = 894.2222
$ws.Range("K3").Value = 894.2222
$ws.Range("M3").Value = -780.2222
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = $null
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = $null
$ws.Range("H107").Value = 1793.2858
$ws.Range("I107").Value = 1321.1818
$ws.Range("J107").Value = 2312.6
$ws.Range("K107").Value = 1321.1818
$ws.Range("L107").Value = 2312.6
$ws.Range("M107").Value = 598.8181999999999
$ws.Range("N107").Value = -6152.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 125001300
$ws.Range("I16").Value = 166667900
$ws.Range("K16").Value = 166667900
$ws.Range("M16").Value = -166667613
$ws.Range("H58").Value = 1096.6666
$ws.Range("I58").Value = 1186.2858
$ws.Range("J58").Value = 971.2
$ws.Range("K58").Value = 1186.2858
$ws.Range("L58").Value = 971.2
$ws.Range("M58").Value = -983.2858000000001
$ws.Range("N58").Value = -1377.2
$ws.Range("H86").Value = 4779185
$ws.Range("I86").Value = 11113913
$ws.Range("J86").Value = 28139.25
$ws.Range("K86").Value = 11113913
$ws.Range("L86").Value = 28139.25
$ws.Range("M86").Value = -11112790
$ws.Range("N86").Value = -30385.25
$ws.Range("H89").Value = 4779185
$ws.Range("I89").Value = 11113913
$ws.Range("J89").Value = 28139.25
$ws.Range("K89").Value = 55569565
$ws.Range("L89").Value = 140696.25
$ws.Range("M89").Value = -55563949
$ws.Range("N89").Value = -151928.25
$ws.Range("H113").Value = 125001300
$ws.Range("I113").Value = 166667900
$ws.Range("K113").Value = 166667900
$ws.Range("M113").Value = -166665730
$ws.Range("H136").Value = 1096.6666
$ws.Range("I136").Value = 1186.2858
$ws.Range("J136").Value = 971.2
$ws.Range("K136").Value = 3558.8574
$ws.Range("L136").Value = 2913.6
$ws.Range("M136").Value = -1008.8574
$ws.Range("N136").Value = -8013.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 556.3871
$ws.Range("I5").Value = 566.2632
$ws.Range("J5").Value = 540.75
$ws.Range("K5").Value = 1698.7896
$ws.Range("L5").Value = 1622.25
$ws.Range("M5").Value = -1586.7896
$ws.Range("N5").Value = -1846.25
$ws.Range("H33").Value = 347.125
$ws.Range("I33").Value = 250.5
$ws.Range("J33").Value = 379.33334
$ws.Range("K33").Value = 1503
$ws.Range("L33").Value = 2276.00004
$ws.Range("M33").Value = -1220
$ws.Range("N33").Value = -2842.00004
$ws.Range("H47").Value = 219.8
$ws.Range("I47").Value = 219.8
$ws.Range("K47").Value = 659.4000000000001
$ws.Range("M47").Value = -228.4000000000001
$ws.Range("H106").Value = 2786
$ws.Range("J106").Value = 2786
$ws.Range("L106").Value = 8358
$ws.Range("N106").Value = -10250
$ws.Range("H131").Value = 37038710
$ws.Range("I131").Value = 142857440
$ws.Range("K131").Value = 428572320
$ws.Range("M131").Value = -428567280
$ws.Range("H135").Value = 556.3871
$ws.Range("I135").Value = 566.2632
$ws.Range("J135").Value = 540.75
$ws.Range("K135").Value = 5096.3688
$ws.Range("L135").Value = 4866.75
$ws.Range("M135").Value = -2561.3688
$ws.Range("N135").Value = -9936.75
$ws.Range("H137").Value = 4511
$ws.Range("I137").Value = 2000
$ws.Range("J137").Value = 5766.5
$ws.Range("K137").Value = 6000
$ws.Range("L137").Value = 17299.5
$ws.Range("M137").Value = -900
$ws.Range("N137").Value = -27499.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2347.8572
$ws.Range("I113").Value = 1417.4286
$ws.Range("K113").Value = 1417.4286
$ws.Range("M113").Value = 752.5714
$ws.Range("H122").Value = 2516.6667
$ws.Range("I122").Value = 2233.3333
$ws.Range("K122").Value = 6699.999899999999
$ws.Range("M122").Value = -4249.999899999999
$ws.Range("H132").Value = 2670.2778
$ws.Range("I132").Value = 2276.2144
$ws.Range("K132").Value = 6828.6432
$ws.Range("M132").Value = -4298.6432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3050.4285
$ws.Range("I40").Value = 2750.6
$ws.Range("K40").Value = 2750.6
$ws.Range("M40").Value = -2614.6
$ws.Range("H61").Value = 1291.1
$ws.Range("I61").Value = 1338.25
$ws.Range("K61").Value = 1338.25
$ws.Range("M61").Value = -1136.25
$ws.Range("H113").Value = 1291.1
$ws.Range("I113").Value = 1338.25
$ws.Range("K113").Value = 1338.25
$ws.Range("M113").Value = 831.75
$ws.Range("H132").Value = 74200.79
$ws.Range("I132").Value = 3750.75
$ws.Range("J132").Value = 102380.8
$ws.Range("K132").Value = 11252.25
$ws.Range("L132").Value = 307142.4
$ws.Range("M132").Value = -8722.25
$ws.Range("N132").Value = -312202.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 627
$ws.Range("I136").Value = 406.45456
$ws.Range("J136").Value = 973.5714
$ws.Range("K136").Value = 1219.36368
$ws.Range("L136").Value = 2920.7142
$ws.Range("M136").Value = 1330.63632
$ws.Range("N136").Value = -8020.7142
